$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PageContent")

# Start from a clean slate for the data region so no stray values survive
# from the previous (smaller) table.
$ws.Range("A1:G40").ClearContents()

# ---------------------------------------------------------------------------
# Rewrite the PageContent data table. The previous "ContentType" (Index /
# ContentTypeDisplayName / Value) rows for the Students Headquarters faculty
# are replaced by a full "mediaDescriptors" style seeder covering photo,
# name, title, education, experience, subjects, office, honors, phone and
# email fields for both full-time and part-time teachers.
# ---------------------------------------------------------------------------

# Header row
$ws.Range("A1").Value = "PageSlug"
$ws.Range("B1").Value = "ContentDefinitionDisplayName"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Value"
$ws.Range("E1").Value = "Sequence"

$slug = "StudentHeadquarters"
$fullTime = "專任教師"
$partTime = "兼任教師"

# Data rows: Name, Value (optional), Sequence
$rows = @(
    @("照片", $null, 0),
    @("姓名", "黃興祿", 0),
    @("職稱", "教授、系主任", 0),
    @("學歷", "中山大學材料科學博士", 0),
    @("經歷", $null, 0),
    @("教授科目", $null, 0),
    @("研究室", $null, 0),
    @("專長領域", "半導體材料、微結構分析", 0),
    @("電話", "(07)7466641", 0),
    @("E-mail", "hlhuang8423@gmail.com", 0),
    @("榮譽", $null, 0),
    @("照片", $null, 1),
    @("姓名", "洪篤傑", 1),
    @("職稱", "教授", 1),
    @("學歷", "中山大學機械工程博士", 1),
    @("經歷", $null, 1),
    @("教授科目", $null, 1),
    @("研究室", $null, 1),
    @("專長領域", "熱流能源、液動拋光法", 1),
    @("電話", "(07)7466641", 1),
    @("E-mail", "evachristrich@gmail.com", 1),
    @("榮譽", $null, 1)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $slug
    $ws.Cells.Item($r, 2).Value = $fullTime
    $ws.Cells.Item($r, 3).Value = $row[0]
    if ($null -ne $row[1]) {
        $ws.Cells.Item($r, 4).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 4).ClearContents()
    }
    $ws.Cells.Item($r, 5).Value = $row[2]
    $r = $r + 1
}

$partTimeRows = @(
    @("照片", $null),
    @("職稱", "副教授"),
    @("姓名", "洪兆宇"),
    @("連絡電話", "(07)7466641"),
    @("E-mail", $null)
)

foreach ($row in $partTimeRows) {
    $ws.Cells.Item($r, 1).Value = $slug
    $ws.Cells.Item($r, 2).Value = $partTime
    $ws.Cells.Item($r, 3).Value = $row[0]
    if ($null -ne $row[1]) {
        $ws.Cells.Item($r, 4).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 4).ClearContents()
    }
    $ws.Cells.Item($r, 5).Value = 0
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Column widths for the new Value / helper columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 27.428571428571427
$ws.Columns.Item(6).ColumnWidth = 19.0
$ws.Columns.Item(7).ColumnWidth = 14.857142857142858

# Active cell ends up past the last populated row, as in the authored sheet
$ws.Range("D29").Select()
